# Apply scheduled market-data refresh to the Leve profit-tracking sheets.
# Each sheet ("ALC", "ARM", "CRP", "CUL", "GSM", "LTW", "WVR") has a table of
# crafting Leve quests with price/profit columns (H:N) populated from an external
# market-data source. This updates the affected rows with refreshed price figures.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 30000
$ws.Range("I21").Value = 30000
$ws.Range("K21").Value = 30000
$ws.Range("M21").Value = -29532

# Row 23
$ws.Range("H23").Value = 30000
$ws.Range("I23").Value = 30000
$ws.Range("K23").Value = 30000
$ws.Range("M23").Value = -29766

# Row 38
$ws.Range("H38").Value = 2253.5557
$ws.Range("J38").Value = 4000
$ws.Range("L38").Value = 12000
$ws.Range("N38").Value = -12744

# Row 43
$ws.Range("H43").Value = 500
$ws.Range("I43").Value = 500
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 500
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -431
$ws.Range("N43").ClearContents()

# Row 61
$ws.Range("H61").Value = 531.25
$ws.Range("I61").Value = 52.5
$ws.Range("J61").Value = 1010
$ws.Range("K61").Value = 157.5
$ws.Range("L61").Value = 3030
$ws.Range("M61").Value = 14.5
$ws.Range("N61").Value = -3374

# Row 116
$ws.Range("H116").Value = 15386499
$ws.Range("I116").Value = 27684758
$ws.Range("J116").Value = 13675
$ws.Range("K116").Value = 27684758
$ws.Range("L116").Value = 13675
$ws.Range("M116").Value = -27681316
$ws.Range("N116").Value = -20559

# Row 123
$ws.Range("H123").Value = 96928.42999999999
$ws.Range("J123").Value = 96928.42999999999
$ws.Range("L123").Value = 96928.42999999999
$ws.Range("N123").Value = -106728.43

# Row 135
$ws.Range("H135").Value = 2747.4
$ws.Range("I135").Value = 2934.25
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 26408.25
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -23873.25
$ws.Range("N135").Value = -23070

# Row 138
$ws.Range("H138").Value = 5884516
$ws.Range("I138").Value = 2075.5334
$ws.Range("J138").Value = 9093120
$ws.Range("K138").Value = 6226.600199999999
$ws.Range("L138").Value = 27279360
$ws.Range("M138").Value = -1086.600199999999
$ws.Range("N138").Value = -27289640

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2269.1482
$ws.Range("I32").Value = 1540.95
$ws.Range("J32").Value = 4349.7144
$ws.Range("K32").Value = 1540.95
$ws.Range("L32").Value = 4349.7144
$ws.Range("M32").Value = -1253.95
$ws.Range("N32").Value = -4923.7144

# Row 74
$ws.Range("H74").Value = 1075.1951
$ws.Range("I74").Value = 1105.6666
$ws.Range("J74").Value = 1016.4286
$ws.Range("K74").Value = 1105.6666
$ws.Range("L74").Value = 1016.4286
$ws.Range("M74").Value = -231.6666
$ws.Range("N74").Value = -2764.4286

# Row 77
$ws.Range("H77").Value = 1075.1951
$ws.Range("I77").Value = 1105.6666
$ws.Range("J77").Value = 1016.4286
$ws.Range("K77").Value = 5528.333000000001
$ws.Range("L77").Value = 5082.143
$ws.Range("M77").Value = -1160.333000000001
$ws.Range("N77").Value = -13818.143

# Row 123
$ws.Range("H123").Value = 39714
$ws.Range("J123").Value = 39714
$ws.Range("L123").Value = 39714
$ws.Range("N123").Value = -49514

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1333.5444
$ws.Range("I31").Value = 1037.8379
$ws.Range("J31").Value = 1539.9811
$ws.Range("K31").Value = 1037.8379
$ws.Range("L31").Value = 1539.9811
$ws.Range("M31").Value = -742.8379
$ws.Range("N31").Value = -2129.9811

# Row 34
$ws.Range("H34").Value = 1333.5444
$ws.Range("I34").Value = 1037.8379
$ws.Range("J34").Value = 1539.9811
$ws.Range("K34").Value = 1037.8379
$ws.Range("L34").Value = 1539.9811
$ws.Range("M34").Value = -835.8379
$ws.Range("N34").Value = -1943.9811

# Row 99
$ws.Range("H99").Value = 6251571
$ws.Range("I99").Value = 8929930
$ws.Range("J99").Value = 2066.6667
$ws.Range("K99").Value = 8929930
$ws.Range("L99").Value = 2066.6667
$ws.Range("M99").Value = -8928432
$ws.Range("N99").Value = -5062.6667

# Row 107
$ws.Range("H107").Value = 507.84616
$ws.Range("J107").Value = 705
$ws.Range("L107").Value = 705
$ws.Range("N107").Value = -4545

# Row 126
$ws.Range("H126").Value = 6251571
$ws.Range("I126").Value = 8929930
$ws.Range("J126").Value = 2066.6667
$ws.Range("K126").Value = 26789790
$ws.Range("L126").Value = 6200.000100000001
$ws.Range("M126").Value = -26787320
$ws.Range("N126").Value = -11140.0001

# Row 132
$ws.Range("H132").Value = 2648.2424
$ws.Range("I132").Value = 2041.45
$ws.Range("J132").Value = 3581.7693
$ws.Range("K132").Value = 6124.35
$ws.Range("L132").Value = 10745.3079
$ws.Range("M132").Value = -3594.35
$ws.Range("N132").Value = -15805.3079

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 8390.486000000001
$ws.Range("J39").Value = 8390.486000000001
$ws.Range("L39").Value = 25171.458
$ws.Range("N39").Value = -25759.458

# Row 68
$ws.Range("H68").Value = 1588.0513
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1588.0513
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 4764.1539
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -6386.1539

# Row 71
$ws.Range("H71").Value = 1588.0513
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1588.0513
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 14292.4617
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -22404.4617

# Row 80
$ws.Range("H80").Value = 1133.3334
$ws.Range("I80").Value = 900
$ws.Range("J80").Value = 1180
$ws.Range("K80").Value = 2700
$ws.Range("L80").Value = 3540
$ws.Range("M80").Value = -1764
$ws.Range("N80").Value = -5412

# Row 83
$ws.Range("H83").Value = 1133.3334
$ws.Range("I83").Value = 900
$ws.Range("J83").Value = 1180
$ws.Range("K83").Value = 8100
$ws.Range("L83").Value = 10620
$ws.Range("M83").Value = -3420
$ws.Range("N83").Value = -19980

# Row 92
$ws.Range("H92").Value = 879.4
$ws.Range("J92").Value = 899.5
$ws.Range("L92").Value = 2698.5
$ws.Range("N92").Value = -5194.5

# Row 131
$ws.Range("H131").Value = 2500.1729
$ws.Range("J131").Value = 2745.1943
$ws.Range("L131").Value = 8235.582900000001
$ws.Range("N131").Value = -18315.5829

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 51
$ws.Range("H51").Value = 58000
$ws.Range("J51").Value = 58000
$ws.Range("L51").Value = 58000
$ws.Range("N51").Value = -59018

# Row 80
$ws.Range("H80").Value = 2120.5715
$ws.Range("I80").Value = 1948.8
$ws.Range("J80").Value = 2550
$ws.Range("K80").Value = 1948.8
$ws.Range("L80").Value = 2550
$ws.Range("M80").Value = -950.8
$ws.Range("N80").Value = -4546

# Row 83
$ws.Range("H83").Value = 2120.5715
$ws.Range("I83").Value = 1948.8
$ws.Range("J83").Value = 2550
$ws.Range("K83").Value = 9744
$ws.Range("L83").Value = 12750
$ws.Range("M83").Value = -4752
$ws.Range("N83").Value = -22734

# Row 102
$ws.Range("H102").Value = 2170.75
$ws.Range("I102").Value = 1388
$ws.Range("J102").Value = 4519
$ws.Range("K102").Value = 1388
$ws.Range("L102").Value = 4519
$ws.Range("M102").Value = 234
$ws.Range("N102").Value = -7763

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 887.0769
$ws.Range("I22").Value = 997.1429000000001
$ws.Range("J22").Value = 758.6667
$ws.Range("K22").Value = 997.1429000000001
$ws.Range("L22").Value = 758.6667
$ws.Range("M22").Value = -702.1429000000001
$ws.Range("N22").Value = -1348.6667

# Row 27
$ws.Range("H27").Value = 887.0769
$ws.Range("I27").Value = 997.1429000000001
$ws.Range("J27").Value = 758.6667
$ws.Range("K27").Value = 997.1429000000001
$ws.Range("L27").Value = 758.6667
$ws.Range("M27").Value = -890.1429000000001
$ws.Range("N27").Value = -972.6667

# Row 46
$ws.Range("H46").Value = 559.8889
$ws.Range("I46").Value = 600
$ws.Range("J46").Value = 548.4286
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 548.4286
$ws.Range("M46").Value = -412
$ws.Range("N46").Value = -924.4286

# Row 55
$ws.Range("H55").Value = 397.2857
$ws.Range("I55").Value = 299.75
$ws.Range("J55").Value = 527.3333
$ws.Range("K55").Value = 299.75
$ws.Range("L55").Value = 527.3333
$ws.Range("M55").Value = -126.75
$ws.Range("N55").Value = -873.3333

# Row 68
$ws.Range("H68").Value = 2036.7273
$ws.Range("I68").Value = 1822.6666
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 1822.6666
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -1073.6666
$ws.Range("N68").Value = -4498

# Row 71
$ws.Range("H71").Value = 2036.7273
$ws.Range("I71").Value = 1822.6666
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 9113.333000000001
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -5369.333000000001
$ws.Range("N71").Value = -22488

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 74
$ws.Range("H74").Value = 14237.2
$ws.Range("J74").Value = 14654.25
$ws.Range("L74").Value = 14654.25
$ws.Range("N74").Value = -16526.25

# Row 77
$ws.Range("H77").Value = 14237.2
$ws.Range("J77").Value = 14654.25
$ws.Range("L77").Value = 43962.75
$ws.Range("N77").Value = -53322.75

# Row 123
$ws.Range("H123").Value = 29904.666
$ws.Range("J123").Value = 29904.666
$ws.Range("L123").Value = 29904.666
$ws.Range("N123").Value = -39704.666

# Row 136
$ws.Range("I136").Value = 23881650
$ws.Range("K136").Value = 71644950
$ws.Range("M136").Value = -71642400
